# "added up to req 7"
#
# Splits the "Sorting list ... all" paragraph so "all" carries its own
# grammar-check proofErr bracket, then appends four new paragraphs:
#   - "Added hyperlinks to each page for easier navigation" (grammar-checked "navigation")
#   - "Css" (spell-checked)
#   - an empty paragraph
#   - "Imported jstl servlet" (spell-checked "jstl")
#
# w:proofErr markers aren't reachable through the normal Range.Text /
# InsertAfter properties, so each paragraph's final OOXML is built as a
# literal <w:p> fragment and dropped in with Range.InsertXML - that is
# the one COM entry point that round-trips raw markup (incl. proofErr)
# verbatim.

$d = $word.ActiveDocument
$W = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function New-WordXml([string]$InnerXml) {
    return '<w:p xmlns:w="' + $W + '">' + $InnerXml + '</w:p>'
}

# 1) Rewrite paragraph 2 in place: "Sorting list to make it easier to view " + [gram]all[/gram]
$p2 = $d.Paragraphs(2)
$p2.Range.InsertXML((New-WordXml (
    '<w:r><w:t xml:space="preserve">Sorting list to make it easier to view </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>all</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>'
)))

# 2) New paragraph: "Added hyperlinks to each page for easier " + [gram]navigation[/gram]
$p2 = $d.Paragraphs(2)
$p2.Range.InsertParagraphAfter()
$p3 = $d.Paragraphs(3)
$p3.Range.InsertXML((New-WordXml (
    '<w:r><w:t xml:space="preserve">Added hyperlinks to each page for easier </w:t></w:r>' +
    '<w:proofErr w:type="gramStart"/>' +
    '<w:r><w:t>navigation</w:t></w:r>' +
    '<w:proofErr w:type="gramEnd"/>'
)))

# 3) New paragraph: [spell]Css[/spell]
$p3 = $d.Paragraphs(3)
$p3.Range.InsertParagraphAfter()
$p4 = $d.Paragraphs(4)
$p4.Range.InsertXML((New-WordXml (
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>Css</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>'
)))

# 4) New, completely empty paragraph
$p4 = $d.Paragraphs(4)
$p4.Range.InsertParagraphAfter()
$p5 = $d.Paragraphs(5)
$p5.Range.InsertXML('<w:p xmlns:w="' + $W + '"/>')

# 5) New paragraph: "Imported " + [spell]jstl[/spell] + " servlet"
$p5 = $d.Paragraphs(5)
$p5.Range.InsertParagraphAfter()
$p6 = $d.Paragraphs(6)
$p6.Range.InsertXML((New-WordXml (
    '<w:r><w:t xml:space="preserve">Imported </w:t></w:r>' +
    '<w:proofErr w:type="spellStart"/>' +
    '<w:r><w:t>jstl</w:t></w:r>' +
    '<w:proofErr w:type="spellEnd"/>' +
    '<w:r><w:t xml:space="preserve"> servlet</w:t></w:r>'
)))
